$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the participant row (row 13: "Emilio Rugerio" second entry) from the
# admin panel. This shifts the following row (Alexis Sharon) up to row 13.
$ws.Rows.Item(13).Delete()
